# edit.ps1
# Refresh the repayment collector figures for the 2025-09-01 .. 2025-09-15
# cycle and mark the worksheet as a duplicated/re-uploaded copy by
# appending " (1)" to its name (mirrors "Add files via upload").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet (re-uploaded copy) -----------------------------
$ws.Name = "repayment_20250901_20250915 (1)"

# --- Updated per-collector repayment figures ------------------------------
# Columns: D=Cycle-days, E=Repayment_amount, F=Pending Amount,
#          G=Pending Amount Recovery, H=Talk_time, I=New_collections,
#          J=Repayment_new_collections, K=New_collection_amount_rate,
#          L=New_collection_count_rate
# Text-look-alike numeric columns (E, F, G, K, L) are written with a
# leading apostrophe so Excel stores them as text (matching the source
# workbook, where these are shared strings like "29,602,130.00"), while
# keeping the cell's number format at General.

$ws.Range("D2").Value = 41
$ws.Range("E2").Value = "'29,602,130.00"
$ws.Range("F2").Value = "'332,666,040.00"
$ws.Range("G2").Value = "'8.90"
$ws.Range("H2").Value = 13.885
$ws.Range("I2").Value = 252
$ws.Range("J2").Value = 15
$ws.Range("K2").Value = "'8.90"
$ws.Range("L2").Value = "'5.95"
$ws.Range("F3").Value = "'325,976,616.00"
$ws.Range("G3").Value = "'9.30"
$ws.Range("H3").Value = 19.15
$ws.Range("I3").Value = 252
$ws.Range("K3").Value = "'2.54"
$ws.Range("L3").Value = "'3.17"
$ws.Range("F4").Value = "'344,627,951.00"
$ws.Range("G4").Value = "'6.55"
$ws.Range("H4").Value = 10.73
$ws.Range("I4").Value = 253
$ws.Range("K4").Value = "'5.37"
$ws.Range("L4").Value = "'6.32"
$ws.Range("F5").Value = "'307,121,314.00"
$ws.Range("G5").Value = "'7.11"
$ws.Range("H5").Value = 19.307
$ws.Range("I5").Value = 222
$ws.Range("K5").Value = "'2.55"
$ws.Range("L5").Value = "'2.70"
$ws.Range("F6").Value = "'316,238,350.00"
$ws.Range("G6").Value = "'7.82"
$ws.Range("H6").Value = 13.501
$ws.Range("I6").Value = 252
$ws.Range("K6").Value = "'2.66"
$ws.Range("L6").Value = "'2.78"
$ws.Range("F7").Value = "'347,149,837.00"
$ws.Range("G7").Value = "'7.62"
$ws.Range("H7").Value = 9.978
$ws.Range("I7").Value = 254
$ws.Range("K7").Value = "'4.32"
$ws.Range("L7").Value = "'4.72"
$ws.Range("F8").Value = "'334,624,168.00"
$ws.Range("G8").Value = "'7.66"
$ws.Range("H8").Value = 19.701
$ws.Range("I8").Value = 249
$ws.Range("K8").Value = "'3.58"
$ws.Range("L8").Value = "'4.42"
$ws.Range("F9").Value = "'356,015,842.00"
$ws.Range("G9").Value = "'7.31"
$ws.Range("H9").Value = 10.809
$ws.Range("I9").Value = 255
$ws.Range("K9").Value = "'5.30"
$ws.Range("L9").Value = "'2.75"
$ws.Range("F10").Value = "'287,278,647.00"
$ws.Range("G10").Value = "'5.39"
$ws.Range("H10").Value = 10.16
$ws.Range("I10").Value = 189
$ws.Range("K10").Value = "'3.47"
$ws.Range("L10").Value = "'3.17"
$ws.Range("F11").Value = "'327,671,563.00"
$ws.Range("G11").Value = "'10.95"
$ws.Range("H11").Value = 11.361
$ws.Range("I11").Value = 250
$ws.Range("K11").Value = "'8.51"
$ws.Range("L11").Value = "'6.80"
$ws.Range("F12").Value = "'315,752,427.00"
$ws.Range("G12").Value = "'7.06"
$ws.Range("H12").Value = 15.418
$ws.Range("I12").Value = 255
$ws.Range("K12").Value = "'6.11"
$ws.Range("L12").Value = "'6.67"
$ws.Range("F13").Value = "'339,780,179.00"
$ws.Range("G13").Value = "'7.11"
$ws.Range("H13").Value = 14.3
$ws.Range("I13").Value = 250
$ws.Range("K13").Value = "'3.62"
$ws.Range("L13").Value = "'3.60"
$ws.Range("D14").Value = 38
$ws.Range("E14").Value = "'25,023,925.00"
$ws.Range("F14").Value = "'324,168,619.00"
$ws.Range("G14").Value = "'7.72"
$ws.Range("H14").Value = 9.126
$ws.Range("I14").Value = 250
$ws.Range("K14").Value = "'5.44"
$ws.Range("L14").Value = "'4.80"
$ws.Range("F15").Value = "'331,674,751.00"
$ws.Range("G15").Value = "'7.15"
$ws.Range("H15").Value = 8.242
$ws.Range("I15").Value = 253
$ws.Range("K15").Value = "'2.99"
$ws.Range("L15").Value = "'3.56"
$ws.Range("F16").Value = "'317,215,129.00"
$ws.Range("G16").Value = "'9.22"
$ws.Range("H16").Value = 8.481
$ws.Range("I16").Value = 254
$ws.Range("K16").Value = "'2.40"
$ws.Range("L16").Value = "'3.94"
$ws.Range("F17").Value = "'319,800,812.00"
$ws.Range("G17").Value = "'9.37"
$ws.Range("H17").Value = 17.812
$ws.Range("I17").Value = 252
$ws.Range("K17").Value = "'2.93"
$ws.Range("L17").Value = "'3.57"
$ws.Range("F18").Value = "'271,676,597.00"
$ws.Range("G18").Value = "'8.14"
$ws.Range("H18").Value = 8.168
$ws.Range("I18").Value = 155
$ws.Range("K18").Value = "'2.83"
$ws.Range("L18").Value = "'1.94"

Write-Output "repayment figures refreshed"
